$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(2000, "falling", -0.2223867177963257, 0.8657857477664948, 0.5839412361383439, -0.0157297793775796, -0.0062613687478005, 0.1360702365636825),
    @(2100, "falling", 0.1144545078277593, 0.9493652880191804, 0.8125896006822588, -0.0059559359215199, -0.0456621758639812, -0.009010262787342),
    @(2200, "falling", 0.001384258270263117, 0.9941152483224868, 0.7851106524467466, 0.0004581489483825, -0.0007635815418325, -0.0609338097274303),
    @(2300, "falling", -0.04791402816772423, 0.9533933401107785, 0.6604413688182829, -0.0006108652451075, 0.0114537235349416, 0.0355829000473022),
    @(2400, "falling", 0.02978670597076406, 0.9465968608856203, 0.6349874883890154, -0.0042760567739605, -0.0050396383740007, -0.0058032199740409),
    @(2500, "falling", 0.007494449615478475, 0.9755845963954928, 0.702362984418869, 0.005192354787141, -0.0007635815418325, -0.0145080499351024),
    @(2600, "falling", 0.001274824142456048, 0.8726469278335574, 0.6991539821028709, -0.0010690141934901, 0.008399397134780801, 0.0061086523346602),
    @(2700, "falling", 0.003990173339843651, 0.9426385164260862, 0.7098991274833678, 0.0004581489483825, 0.0038179077673703, -0.0076358155347406),
    @(2800, "falling", -0.00699806213378908, 0.904570586979389, 0.691669926047325, 0.0007635815418325, 0.0027488935738801, -0.0044287731871008),
    @(2900, "falling", 0.01654338836669934, 0.8796091377735137, 0.693576380610466, -0.00167987938039, -0.0027488935738801, 0.0114537235349416)
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
